# Fussenegger_Task1_RS1_SOSE13 - "Updated german translation of KDD."
#
# The slide "Clustering in Spatial Databases" (slide 2) lists bullet
# points about spatial databases; the line
#   "Wissensfindung in Datenbanken"
# is retranslated to
#   "Wissensentdeckung in Datenbanken"
# by retyping just the first word ("Wissensfindung" -> "Wissensentdeckung"),
# leaving " in Datenbanken" (and the rest of the paragraph - the
# line break plus the "(knowledge discovery in databases)" gloss) as-is.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text.IndexOf("Wissensfindung in Datenbanken") -ge 0) {
                $targetSlide = $s
                $targetShape = $sh
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the 'Wissensfindung in Datenbanken' text on any slide"
}

$tr = $targetShape.TextFrame.TextRange

# Re-type the "Wissensfindung " part of the run in place; PowerPoint
# splits the original run at the edit boundary, so " in Datenbanken"
# keeps its original run properties untouched.
$idx = $tr.Text.IndexOf("Wissensfindung in Datenbanken")
$old = $tr.Characters($idx + 1, "Wissensfindung ".Length)
$old.Text = "Wissensentdeckung "
